# Workbook: "Hortaliza, Terminal La Palmera de La Serena - Brócoli"
# This weekly update inserts a new pair of records (Primera/Segunda quality)
# at the top of the data block (rows 257-258), pushing all the existing
# records down by two rows (from A1:R382 to A1:R384).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 257-258, shifting rows 257:382 down to 259:384.
$insertRange = $ws.Range("A257:R258")
$insertRange.Insert(-4121)  # xlShiftDown

# Populate the new row 257 (Calidad = Primera)
$ws.Cells.Item(257, 1).Value2 = 8
$ws.Cells.Item(257, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(257, 3).Value2 = "Coquimbo"
$ws.Cells.Item(257, 4).Value2 = 44460
$ws.Cells.Item(257, 5).Value2 = 4
$ws.Cells.Item(257, 6).Value2 = 100112023
$ws.Cells.Item(257, 7).Value2 = "Brócoli"
$ws.Cells.Item(257, 8).Value2 = "Sin especificar"
$ws.Cells.Item(257, 9).Value2 = "Primera"
$ws.Cells.Item(257, 10).Value2 = 2000
$ws.Cells.Item(257, 11).Value2 = 600
$ws.Cells.Item(257, 12).Value2 = 700
$ws.Cells.Item(257, 13).Value2 = 650
$ws.Cells.Item(257, 14).Value2 = "`$/unidad"
$ws.Cells.Item(257, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(257, 16).Value2 = 650
$ws.Cells.Item(257, 17).Value2 = 1
$ws.Cells.Item(257, 18).Value2 = "Hortaliza"

# Populate the new row 258 (Calidad = Segunda)
$ws.Cells.Item(258, 1).Value2 = 8
$ws.Cells.Item(258, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(258, 3).Value2 = "Coquimbo"
$ws.Cells.Item(258, 4).Value2 = 44460
$ws.Cells.Item(258, 5).Value2 = 4
$ws.Cells.Item(258, 6).Value2 = 100112023
$ws.Cells.Item(258, 7).Value2 = "Brócoli"
$ws.Cells.Item(258, 8).Value2 = "Sin especificar"
$ws.Cells.Item(258, 9).Value2 = "Segunda"
$ws.Cells.Item(258, 10).Value2 = 1400
$ws.Cells.Item(258, 11).Value2 = 500
$ws.Cells.Item(258, 12).Value2 = 550
$ws.Cells.Item(258, 13).Value2 = 525
$ws.Cells.Item(258, 14).Value2 = "`$/unidad"
$ws.Cells.Item(258, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(258, 16).Value2 = 525
$ws.Cells.Item(258, 17).Value2 = 1
$ws.Cells.Item(258, 18).Value2 = "Hortaliza"

Write-Host "Done. New dimension rows should span to 384."
